$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.27784
$ws.Range("H2").Value = 30.83352
$ws.Range("I2").Value = 0.230301226653591
$ws.Range("J2").Value = 0.230301226653591
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.020559
$ws.Range("N2").Value = 90.061677
$ws.Range("O2").Value = 0.8829766276144534
$ws.Range("P2").Value = 0.8829766276144534
$ws.Range("Q2").Value = 308.54650211256
$ws.Range("R2").Value = 2776.91851901304
$ws.Range("S2").Value = 0.2033506004460597
$ws.Range("T2").Value = 0.2033506004460596

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.27784
$ws.Range("H3").Value = 30.83352
$ws.Range("I3").Value = 0.230301226653591
$ws.Range("J3").Value = 0.230301226653591
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.426936666666667
$ws.Range("N3").Value = 4.28081
$ws.Range("O3").Value = 0.04196962907162197
$ws.Range("P3").Value = 0.04196962907162197
$ws.Range("Q3").Value = 14.66582675013333
$ws.Range("R3").Value = 131.9924407512
$ws.Range("S3").Value = 0.009665657057390753
$ws.Range("T3").Value = 0.009665657057390753

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.27784
$ws.Range("H4").Value = 30.83352
$ws.Range("I4").Value = 0.230301226653591
$ws.Range("J4").Value = 0.230301226653591
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.067853
$ws.Range("N4").Value = 3.203559
$ws.Range("O4").Value = 0.03140811737476231
$ws.Range("P4").Value = 0.0314081173747623
$ws.Range("Q4").Value = 10.97522227752
$ws.Range("R4").Value = 98.77700049768001
$ws.Range("S4").Value = 0.007233327958287725
$ws.Range("T4").Value = 0.007233327958287722

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.27784
$ws.Range("H5").Value = 30.83352
$ws.Range("I5").Value = 0.230301226653591
$ws.Range("J5").Value = 0.230301226653591
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.483919333333333
$ws.Range("N5").Value = 4.451758
$ws.Range("O5").Value = 0.04364562593916237
$ws.Range("P5").Value = 0.04364562593916237
$ws.Range("Q5").Value = 15.25148548090667
$ws.Range("R5").Value = 137.26336932816
$ws.Range("S5").Value = 0.01005164119185288
$ws.Range("T5").Value = 0.01005164119185288

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.56812733333333
$ws.Range("H6").Value = 43.704382
$ws.Range("I6").Value = 0.3264360600001921
$ws.Range("J6").Value = 0.326436060000192
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.020559
$ws.Range("N6").Value = 90.061677
$ws.Range("O6").Value = 0.8829766276144534
$ws.Range("P6").Value = 0.8829766276144534
$ws.Range("Q6").Value = 437.343326129846
$ws.Range("R6").Value = 3936.089935168614
$ws.Range("S6").Value = 0.288235411390719
$ws.Range("T6").Value = 0.288235411390719

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.56812733333333
$ws.Range("H7").Value = 43.704382
$ws.Range("I7").Value = 0.3264360600001921
$ws.Range("J7").Value = 0.326436060000192
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.426936666666667
$ws.Range("N7").Value = 4.28081
$ws.Range("O7").Value = 0.04196962907162197
$ws.Range("P7").Value = 0.04196962907162197
$ws.Range("Q7").Value = 20.78779505660222
$ws.Range("R7").Value = 187.09015550942
$ws.Range("S7").Value = 0.0137004003538098
$ws.Range("T7").Value = 0.0137004003538098

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.56812733333333
$ws.Range("H8").Value = 43.704382
$ws.Range("I8").Value = 0.3264360600001921
$ws.Range("J8").Value = 0.326436060000192
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.067853
$ws.Range("N8").Value = 3.203559
$ws.Range("O8").Value = 0.03140811737476231
$ws.Range("P8").Value = 0.0314081173747623
$ws.Range("Q8").Value = 15.556618477282
$ws.Range("R8").Value = 140.009566295538
$ws.Range("S8").Value = 0.01025274208784098
$ws.Range("T8").Value = 0.01025274208784098

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.56812733333333
$ws.Range("H9").Value = 43.704382
$ws.Range("I9").Value = 0.3264360600001921
$ws.Range("J9").Value = 0.326436060000192
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.483919333333333
$ws.Range("N9").Value = 4.451758
$ws.Range("O9").Value = 0.04364562593916237
$ws.Range("P9").Value = 0.04364562593916237
$ws.Range("Q9").Value = 21.61792580039511
$ws.Range("R9").Value = 194.561332203556
$ws.Range("S9").Value = 0.01424750616782235
$ws.Range("T9").Value = 0.01424750616782235

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.650137
$ws.Range("H10").Value = 7.950411
$ws.Range("I10").Value = 0.05938308067649115
$ws.Range("J10").Value = 0.05938308067649114
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.020559
$ws.Range("N10").Value = 90.061677
$ws.Range("O10").Value = 0.8829766276144534
$ws.Range("P10").Value = 0.8829766276144534
$ws.Range("Q10").Value = 79.55859416658301
$ws.Range("R10").Value = 716.0273474992471
$ws.Range("S10").Value = 0.05243387231308517
$ws.Range("T10").Value = 0.05243387231308516

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.650137
$ws.Range("H11").Value = 7.950411
$ws.Range("I11").Value = 0.05938308067649115
$ws.Range("J11").Value = 0.05938308067649114
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.426936666666667
$ws.Range("N11").Value = 4.28081
$ws.Range("O11").Value = 0.04196962907162197
$ws.Range("P11").Value = 0.04196962907162197
$ws.Range("Q11").Value = 3.78157765699
$ws.Range("R11").Value = 34.03419891291
$ws.Range("S11").Value = 0.002492285869122536
$ws.Range("T11").Value = 0.002492285869122536

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.650137
$ws.Range("H12").Value = 7.950411
$ws.Range("I12").Value = 0.05938308067649115
$ws.Range("J12").Value = 0.05938308067649114
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.067853
$ws.Range("N12").Value = 3.203559
$ws.Range("O12").Value = 0.03140811737476231
$ws.Range("P12").Value = 0.0314081173747623
$ws.Range("Q12").Value = 2.829956745861
$ws.Range("R12").Value = 25.469610712749
$ws.Range("S12").Value = 0.001865110767962214
$ws.Range("T12").Value = 0.001865110767962213

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.650137
$ws.Range("H13").Value = 7.950411
$ws.Range("I13").Value = 0.05938308067649115
$ws.Range("J13").Value = 0.05938308067649114
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.483919333333333
$ws.Range("N13").Value = 4.451758
$ws.Range("O13").Value = 0.04364562593916237
$ws.Range("P13").Value = 0.04364562593916237
$ws.Range("Q13").Value = 3.932589530282
$ws.Range("R13").Value = 35.393305772538
$ws.Range("S13").Value = 0.002591811726321234
$ws.Range("T13").Value = 0.002591811726321233

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.13170833333333
$ws.Range("H14").Value = 51.395125
$ws.Range("I14").Value = 0.3838796326697257
$ws.Range("J14").Value = 0.3838796326697257
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.020559
$ws.Range("N14").Value = 90.061677
$ws.Range("O14").Value = 0.8829766276144534
$ws.Range("P14").Value = 0.8829766276144534
$ws.Range("Q14").Value = 514.303460791625
$ws.Range("R14").Value = 4628.731147124625
$ws.Range("S14").Value = 0.3389567434645896
$ws.Range("T14").Value = 0.3389567434645895

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.13170833333333
$ws.Range("H15").Value = 51.395125
$ws.Range("I15").Value = 0.3838796326697257
$ws.Range("J15").Value = 0.3838796326697257
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.426936666666667
$ws.Range("N15").Value = 4.28081
$ws.Range("O15").Value = 0.04196962907162197
$ws.Range("P15").Value = 0.04196962907162197
$ws.Range("Q15").Value = 24.44586278347222
$ws.Range("R15").Value = 220.01276505125
$ws.Range("S15").Value = 0.01611128579129888
$ws.Range("T15").Value = 0.01611128579129888

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.13170833333333
$ws.Range("H16").Value = 51.395125
$ws.Range("I16").Value = 0.3838796326697257
$ws.Range("J16").Value = 0.3838796326697257
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.067853
$ws.Range("N16").Value = 3.203559
$ws.Range("O16").Value = 0.03140811737476231
$ws.Range("P16").Value = 0.0314081173747623
$ws.Range("Q16").Value = 18.294146138875
$ws.Range("R16").Value = 164.647315249875
$ws.Range("S16").Value = 0.01205693656067139
$ws.Range("T16").Value = 0.01205693656067138

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.13170833333333
$ws.Range("H17").Value = 51.395125
$ws.Range("I17").Value = 0.3838796326697257
$ws.Range("J17").Value = 0.3838796326697257
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.483919333333333
$ws.Range("N17").Value = 4.451758
$ws.Range("O17").Value = 0.04364562593916237
$ws.Range("P17").Value = 0.04364562593916237
$ws.Range("Q17").Value = 25.42207320886111
$ws.Range("R17").Value = 194.561332203556
$ws.Range("S17").Value = 0.0167546668531659
$ws.Range("T17").Value = 0.0167546668531659

Write-Output "Updated cells for rows 2-17"
